# Updated cryptos list on Mon Sep 11 23:38:36 UTC 2023 with GitHub Actions
#
# Refreshes the crypto price/volume table: column D (Price) and column E
# (Volume(1h)) are updated with newly scraped values for most rows, and
# rows 40-41 swap their coin (ARBITRUM / MXToken) content entirely.
#
# Price values are stored as text in the workbook (they use dotted
# thousand separators, e.g. "25.239.95", or need to preserve trailing
# zeros, e.g. "0.0780"). When a new price string looks like a plain
# number, Excel would otherwise silently convert it to a numeric value
# (losing formatting / precision), so for those cells we force the
# cell's number format to Text ("@") before assigning the value, then
# restore the default "Normal" style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.239.95"
$ws.Range("E2").Value = "  -2.90%  "
$ws.Range("D3").Value = "1.554.79"
$ws.Range("E3").Value = "  -4.46%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.29%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.479"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0609"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.243"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0780"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("E12").Value = "  -4.49%  "
$ws.Range("D13").Value = "1.550.99"
$ws.Range("E13").Value = "  -4.89%  "
$ws.Range("E14").Value = "  -4.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.505"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.06%  "
$ws.Range("D16").Value = "25.242.65"
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "58.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.18%  "
$ws.Range("D18").Value = "0.0₃0707"
$ws.Range("E18").Value = "  -4.63%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "185.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.53%  "
$ws.Range("E21").Value = "  -3.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.130"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.77%  "
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("E27").Value = "  -4.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "14.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("E29").Value = "  -4.73%  "
$ws.Range("E31").Value = "  -3.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.32%  "
$ws.Range("E33").Value = "  -4.98%  "
$ws.Range("E34").Value = "  -3.21%  "
$ws.Range("E35").Value = "  -3.81%  "
$ws.Range("D36").Value = "1.083.65"
$ws.Range("E36").Value = "  -3.56%  "
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0149"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.27%  "
$ws.Range("E39").Value = "  -4.85%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.47%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.763"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.21%  "
$ws.Range("E42").Value = "  +5.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "92.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("D45").Value = "1.685.03"
$ws.Range("E45").Value = "  -4.45%  "
$ws.Range("E46").Value = "  -2.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E49").Value = "  -3.88%  "
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("E51").Value = "  -2.20%  "
